# Auto-generated Excel COM-interop script
# Applies the "cryptos list" update commit (Fri Jul 28 19:23:00 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.314.45"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.876.20"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'242.39"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.3111"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'0.07734"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "'25.09"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "'0.08440"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").Value = "1.881.83"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "'5.216"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "'0.7112"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "'91.38"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "29.312.94"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "'0.000008289"
$ws.Range("E17").Value = "  +6.40%  "
$ws.Range("D18").Value = "'5.990"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("D19").Value = "'242.75"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "2.127.76"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'7.815"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'0.1620"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").Value = "'163.25"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").Value = "'9.022"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").Value = "'18.52"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "'4.421"
$ws.Range("D31").Value = "'4.334"
$ws.Range("E31").Value = "  +6.09%  "
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("D33").Value = "'0.05259"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7481"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.173"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "'2.683"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "'0.01860"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "'2.725"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").Value = "1.161.76"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").Value = "'6.368"
$ws.Range("E41").Value = "  +4.24%  "
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "'0.8894"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "'106.50"
$ws.Range("E44").Value = "  +4.73%  "
$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "2.024.45"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").Value = "'1.808"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").Value = "'0.5192"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").Value = "'0.00000000120"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").Value = "'9.393"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").Value = "'0.4300"
$ws.Range("E51").Value = "  +1.49%  "
